$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Columns that look numeric/date-like need to be forced to Text so Excel's
# auto-conversion doesn't strip the leading zero / turn the date into a
# serial number. Columns with non-numeric-looking text (region, description,
# evaluation amount with a space, name) and the genuinely numeric "Item
# Number" column don't need this.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01120015455"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "112"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "385"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02.12.2025"

$ws.Cells.Item($row, 5).Value = "Санкт-Петербург"

$ws.Cells.Item($row, 6).Value = 1

$ws.Cells.Item($row, 7).Value = "Тест"

$ws.Cells.Item($row, 8).Value = "1 000"

$ws.Cells.Item($row, 9).Value = "Oleg"
